# Updates cryptos list price (D) and 1h-volume-change (E) columns for Sheet1.
# A handful of Price values are digit-strings that Excel would otherwise
# auto-coerce into numbers (e.g. "1.000" -> 1), so those are written with a
# leading apostrophe to force literal text, matching the source inlineStr cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.179.84"
$ws.Range("E2").Value = "  -0.58%  "
$ws.Range("D3").Value = "1.855.26"
$ws.Range("D4").Value = "0.9999"
$ws.Range("E4").Value = "  -0.20%  "
$ws.Range("D5").Value = "237.79"
$ws.Range("E5").Value = "  -0.52%  "
$ws.Range("D6").Value = "0.6889"
$ws.Range("E6").Value = "  -1.38%  "
$ws.Range("D7").Value = "'1.000"
$ws.Range("E7").Value = "  -0.23%  "
$ws.Range("D8").Value = "0.07769"
$ws.Range("E8").Value = "  +5.14%  "
$ws.Range("D9").Value = "'0.3050"
$ws.Range("E9").Value = "  -0.55%  "
$ws.Range("D10").Value = "23.21"
$ws.Range("E10").Value = "  -2.01%  "
$ws.Range("D11").Value = "0.08068"
$ws.Range("E11").Value = "  -0.52%  "
$ws.Range("D12").Value = "1.861.19"
$ws.Range("E12").Value = "  -1.46%  "
$ws.Range("D13").Value = "0.7209"
$ws.Range("E13").Value = "  -0.53%  "
$ws.Range("D14").Value = "5.184"
$ws.Range("E14").Value = "  -0.52%  "
$ws.Range("D15").Value = "'89.30"
$ws.Range("E15").Value = "  -0.29%  "
$ws.Range("D16").Value = "29.181.11"
$ws.Range("E16").Value = "  -1.08%  "
$ws.Range("D17").Value = "5.733"
$ws.Range("E17").Value = "  -2.85%  "
$ws.Range("D18").Value = "'0.000007802"
$ws.Range("E18").Value = "  +1.10%  "
$ws.Range("D19").Value = "13.24"
$ws.Range("E19").Value = "  +1.08%  "
$ws.Range("D20").Value = "234.68"
$ws.Range("D21").Value = "0.9991"
$ws.Range("E21").Value = "  -0.31%  "
$ws.Range("D22").Value = "2.112.29"
$ws.Range("E22").Value = "  -1.07%  "
$ws.Range("D23").Value = "'1.000"
$ws.Range("E23").Value = "  -0.21%  "
$ws.Range("D24").Value = "7.462"
$ws.Range("E24").Value = "  -2.01%  "
$ws.Range("D25").Value = "162.12"
$ws.Range("E25").Value = "  +0.61%  "
$ws.Range("D26").Value = "8.962"
$ws.Range("E26").Value = "  -0.71%  "
$ws.Range("D27").Value = "0.1424"
$ws.Range("D28").Value = "18.02"
$ws.Range("E28").Value = "  -0.20%  "
$ws.Range("E29").Value = "  +0.47%  "
$ws.Range("D30").Value = "1.398"
$ws.Range("E30").Value = "  +0.79%  "
$ws.Range("D31").Value = "4.516"
$ws.Range("E31").Value = "  +2.40%  "
$ws.Range("D32").Value = "1.481"
$ws.Range("E32").Value = "  -1.70%  "
$ws.Range("D33").Value = "4.006"
$ws.Range("E33").Value = "  -1.30%  "
$ws.Range("D34").Value = "0.05203"
$ws.Range("E34").Value = "  -1.27%  "
$ws.Range("D35").Value = "1.183"
$ws.Range("E35").Value = "  -1.33%  "
$ws.Range("D36").Value = "0.7026"
$ws.Range("E36").Value = "  -2.70%  "
$ws.Range("D37").Value = "1.007"
$ws.Range("E37").Value = "  +0.26%  "
$ws.Range("D38").Value = "2.676"
$ws.Range("E38").Value = "  -0.23%  "
$ws.Range("D39").Value = "0.01844"
$ws.Range("E39").Value = "  -1.11%  "
$ws.Range("D40").Value = "2.679"
$ws.Range("D41").Value = "0.9397"
$ws.Range("E41").Value = "  +7.44%  "
$ws.Range("D42").Value = "1.090.72"
$ws.Range("E42").Value = "  +5.18%  "
$ws.Range("D43").Value = "'5.980"
$ws.Range("E43").Value = "  +0.90%  "
$ws.Range("D44").Value = "0.4282"
$ws.Range("E44").Value = "  -0.64%  "
$ws.Range("D45").Value = "70.41"
$ws.Range("E45").Value = "  +0.90%  "
$ws.Range("D46").Value = "0.9999"
$ws.Range("E46").Value = "  -0.19%  "
$ws.Range("D47").Value = "102.34"
$ws.Range("E47").Value = "  +0.01%  "
$ws.Range("D48").Value = "1.793"
$ws.Range("E48").Value = "  +1.86%  "
$ws.Range("D49").Value = "2.008.54"
$ws.Range("E49").Value = "  -1.06%  "
$ws.Range("D50").Value = "9.149"
$ws.Range("E50").Value = "  -1.09%  "
$ws.Range("D51").Value = "7.001"
$ws.Range("E51").Value = "  -3.55%  "
